# Scheduled-runner refresh: re-sync cached FFXIV Market Board price snapshots
# (currentAveragePrice[/NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ] -- cols H:N)
# for the leve rows whose underlying item prices moved since the last run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64: Forged from the Void (Void Glue)
$ws.Range("H64").Value = 2603.4482
$ws.Range("I64").Value = 2515
$ws.Range("J64").Value = 2800
$ws.Range("K64").Value = 2515
$ws.Range("L64").Value = 2800
$ws.Range("M64").Value = -2267
$ws.Range("N64").Value = -3296

# Row 67: Dodging the Draft (L) (Void Glue)
$ws.Range("H67").Value = 2603.4482
$ws.Range("I67").Value = 2515
$ws.Range("J67").Value = 2800
$ws.Range("K67").Value = 2515
$ws.Range("L67").Value = 2800
$ws.Range("M67").Value = -1657
$ws.Range("N67").Value = -4516

# Row 74: Adhesive of Antipathy (Wing Glue)
$ws.Range("H74").Value = 4003
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# Row 77: It's Gonna Grow Back (L) (Wing Glue)
$ws.Range("H77").Value = 4003
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust (Steel Ingot)
$ws.Range("H32").Value = 13162403
$ws.Range("I32").Value = 3735.767
$ws.Range("J32").Value = 333356640
$ws.Range("K32").Value = 3735.767
$ws.Range("L32").Value = 333356640
$ws.Range("M32").Value = -3448.767
$ws.Range("N32").Value = -333357214

# Row 97: Ore for Me (High Steel Ingot)
$ws.Range("H97").Value = 1715.7142
$ws.Range("I97").Value = 2047
$ws.Range("J97").Value = 1384.4286
$ws.Range("K97").Value = 2047
$ws.Range("L97").Value = 1384.4286
$ws.Range("M97").Value = -1551
$ws.Range("N97").Value = -2376.4286

# Row 132: Don't Bore Me, Ore Me (Mountain Chromite Ingot)
$ws.Range("H132").Value = 1051542.8
$ws.Range("I132").Value = 1053.2683
$ws.Range("J132").Value = 3922880.5
$ws.Range("K132").Value = 3159.8049
$ws.Range("L132").Value = 11768641.5
$ws.Range("M132").Value = -629.8049000000001
$ws.Range("N132").Value = -11773701.5

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt (Iron Ingot)
$ws.Range("H20").Value = 2639.24
$ws.Range("I20").Value = 2416
$ws.Range("J20").Value = 3113.625
$ws.Range("K20").Value = 2416
$ws.Range("L20").Value = 3113.625
$ws.Range("M20").Value = -2169
$ws.Range("N20").Value = -3607.625

# Row 86: Through Thick and Thin (Adamantite Nugget)
$ws.Range("H86").Value = 803444.4399999999
$ws.Range("I86").Value = 1282.1904
$ws.Range("J86").Value = 2909120.2
$ws.Range("K86").Value = 1282.1904
$ws.Range("L86").Value = 2909120.2
$ws.Range("M86").Value = -159.1904
$ws.Range("N86").Value = -2911366.2

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) (Adamantite Nugget)
$ws.Range("H89").Value = 803444.4399999999
$ws.Range("I89").Value = 1282.1904
$ws.Range("J89").Value = 2909120.2
$ws.Range("K89").Value = 6410.951999999999
$ws.Range("L89").Value = 14545601
$ws.Range("M89").Value = -794.9519999999993
$ws.Range("N89").Value = -14556833

# Row 94: High Steal (High Steel Nugget)
$ws.Range("H94").Value = 817.8421
$ws.Range("I94").Value = 648.8889
$ws.Range("J94").Value = 969.9
$ws.Range("K94").Value = 648.8889
$ws.Range("L94").Value = 969.9
$ws.Range("M94").Value = -197.8889
$ws.Range("N94").Value = -1871.9

# Row 99: Meddle in Metal (Oroshigane Ingot)
$ws.Range("H99").Value = 2154.5454
$ws.Range("I99").Value = 800
$ws.Range("J99").Value = 3283.3333
$ws.Range("K99").Value = 800
$ws.Range("L99").Value = 3283.3333
$ws.Range("M99").Value = 698
$ws.Range("N99").Value = -6279.3333

# Row 134: Ruthenium Supremium (Ruthenium Ingot)
$ws.Range("H134").Value = 5788.1904
$ws.Range("I134").Value = 891.2941
$ws.Range("J134").Value = 26600
$ws.Range("K134").Value = 2673.8823
$ws.Range("L134").Value = 79800
$ws.Range("M134").Value = -138.8822999999998
$ws.Range("N134").Value = -84870

$ws = $wb.Worksheets.Item("CRP")
# Row 62: Splinter in the Sewers (Cedar Lumber)
$ws.Range("H62").Value = 4311.2856
$ws.Range("I62").Value = 2599
$ws.Range("J62").Value = 4596.6665
$ws.Range("K62").Value = 2599
$ws.Range("L62").Value = 4596.6665
$ws.Range("M62").Value = -1975
$ws.Range("N62").Value = -5844.6665

# Row 65: The Lumber of Their Discontent (L) (Cedar Lumber)
$ws.Range("H65").Value = 4311.2856
$ws.Range("I65").Value = 2599
$ws.Range("J65").Value = 4596.6665
$ws.Range("K65").Value = 12995
$ws.Range("L65").Value = 22983.3325
$ws.Range("M65").Value = -9875
$ws.Range("N65").Value = -29223.3325

# Row 105: Zelkova, My Love (Zelkova Lumber)
$ws.Range("H105").Value = 7388.2666
$ws.Range("I105").Value = 7994.923
$ws.Range("K105").Value = 7994.923
$ws.Range("M105").Value = -6247.923

$ws = $wb.Worksheets.Item("CUL")
# Row 3: Trout Fishing in Limsa (Grilled Trout)
$ws.Range("H3").Value = 5599.5
$ws.Range("I3").Value = 5919.4
$ws.Range("J3").Value = 4000
$ws.Range("K3").Value = 17758.2
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = -17646.2
$ws.Range("N3").Value = -12224

# Row 55: Pagan Pastries (Pastry Fish)
$ws.Range("H55").Value = 14309.777
$ws.Range("I55").Value = 800
$ws.Range("J55").Value = 15998.5
$ws.Range("K55").Value = 2400
$ws.Range("L55").Value = 47995.5
$ws.Range("M55").Value = -2223
$ws.Range("N55").Value = -48349.5

# Row 113: Can't Eat Just One (Night Vinegar)
$ws.Range("H113").Value = 16667158
$ws.Range("I113").Value = 477.84616
$ws.Range("J113").Value = 33333838
$ws.Range("K113").Value = 1433.53848
$ws.Range("L113").Value = 100001514
$ws.Range("M113").Value = 736.4615200000001
$ws.Range("N113").Value = -100005854

# Row 131: The Mountain Steeped (Tsai tou Vounou)
$ws.Range("H131").Value = 764.2
$ws.Range("J131").Value = 792.32965
$ws.Range("L131").Value = 2376.98895
$ws.Range("N131").Value = -12456.98895

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell (Hardsilver Ingot)
$ws.Range("H80").Value = 5557427.5
$ws.Range("I80").Value = 1976.9231
$ws.Range("J80").Value = 20001600
$ws.Range("K80").Value = 1976.9231
$ws.Range("L80").Value = 20001600
$ws.Range("M80").Value = -978.9231
$ws.Range("N80").Value = -20003596

# Row 83: With a Noise That Reaches Heaven (L) (Hardsilver Ingot)
$ws.Range("H83").Value = 5557427.5
$ws.Range("I83").Value = 1976.9231
$ws.Range("J83").Value = 20001600
$ws.Range("K83").Value = 9884.6155
$ws.Range("L83").Value = 100008000
$ws.Range("M83").Value = -4892.6155
$ws.Range("N83").Value = -100017984

# Row 132: On Board for Lar (Lar Ingot)
$ws.Range("H132").Value = 6325.778
$ws.Range("I132").Value = 2572.923
$ws.Range("J132").Value = 9810.571
$ws.Range("K132").Value = 7718.768999999999
$ws.Range("L132").Value = 29431.713
$ws.Range("M132").Value = -5188.768999999999
$ws.Range("N132").Value = -34491.713

$ws = $wb.Worksheets.Item("LTW")
# Row 45: Soft Shoe Shuffle (Peisteskin Crakows)
$ws.Range("H45").Value = 4250
$ws.Range("I45").Value = 666.6667
$ws.Range("J45").Value = 15000
$ws.Range("K45").Value = 666.6667
$ws.Range("L45").Value = 15000
$ws.Range("M45").Value = -259.6667
$ws.Range("N45").Value = -15814

# Row 55: It's Not a Job, It's a Calling (Peiste Leather)
$ws.Range("H55").Value = 41670940
$ws.Range("I55").Value = 12593.75
$ws.Range("J55").Value = 62500110
$ws.Range("K55").Value = 12593.75
$ws.Range("L55").Value = 62500110
$ws.Range("M55").Value = -12420.75
$ws.Range("N55").Value = -62500456

# Row 68: You Could Say It's a Moving Target (Wyvern Leather)
$ws.Range("H68").Value = 1435.0555
$ws.Range("I68").Value = 1460.0625
$ws.Range("J68").Value = 1235
$ws.Range("K68").Value = 1460.0625
$ws.Range("L68").Value = 1235
$ws.Range("M68").Value = -711.0625
$ws.Range("N68").Value = -2733

# Row 71: They Call It Bloody Mary (L) (Wyvern Leather)
$ws.Range("H71").Value = 1435.0555
$ws.Range("I71").Value = 1460.0625
$ws.Range("J71").Value = 1235
$ws.Range("K71").Value = 7300.3125
$ws.Range("L71").Value = 6175
$ws.Range("M71").Value = -3556.3125
$ws.Range("N71").Value = -13663

# Row 93: Hide to Go Seek (Gagana Leather)
$ws.Range("H93").Value = 866.1429000000001
$ws.Range("I93").Value = 808.2222
$ws.Range("J93").Value = 970.4
$ws.Range("K93").Value = 808.2222
$ws.Range("L93").Value = 970.4
$ws.Range("M93").Value = 439.7778
$ws.Range("N93").Value = -3466.4

$ws = $wb.Worksheets.Item("WVR")
# Row 96: Skills on Display (Ruby Cotton Cloth)
$ws.Range("H96").Value = 2100
$ws.Range("I96").Value = 1520
$ws.Range("J96").Value = 5000
$ws.Range("K96").Value = 1520
$ws.Range("L96").Value = 5000
$ws.Range("M96").Value = -147
$ws.Range("N96").Value = -7746
